# Update "paises.xlsx" COVID stats sheet: refresh the timestamp banner,
# update several countries' case counts, and re-rank Croacia/Zambia now
# that Croacia's total cases have overtaken Zambia's.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Septiembre de 2020 a las 10:29"

# India (row 5)
$ws.Range("B5").Value = 4663930
$ws.Range("C5").Value = 6551
$ws.Range("D5").Value = 3624375
$ws.Range("E5").Value = 962018
$ws.Range("G5").Value = 31
$ws.Range("H5").Value = 77537

# Rusia (row 7)
$ws.Range("B7").Value = 1057362
$ws.Range("C7").Value = 5488
$ws.Range("D7").Value = 873535
$ws.Range("E7").Value = 165343
$ws.Range("G7").Value = 119
$ws.Range("H7").Value = 18484

# Filipinas (row 25)
$ws.Range("B25").Value = 257863
$ws.Range("C25").Value = 4935
$ws.Range("D25").Value = 187116
$ws.Range("E25").Value = 66455
$ws.Range("G25").Value = 186
$ws.Range("H25").Value = 4292

# Barein (row 53)
$ws.Range("E53").Value = 5853
$ws.Range("G53").Value = 2
$ws.Range("H53").Value = 210

# Singapur (row 55)
$ws.Range("B55").Value = 57357
$ws.Range("C55").Value = 42
$ws.Range("E55").Value = 723

# Afganistan (row 66)
$ws.Range("B66").Value = 38641
$ws.Range("C66").Value = 35
$ws.Range("D66").Value = 31234
$ws.Range("E66").Value = 5987

# Croacia overtakes Zambia (Croacia's total cases rise to 13368, above
# Zambia's 13323) so the two swap places: row 89 becomes Croacia with its
# updated numbers, row 90 becomes Zambia with its prior (unchanged) numbers.
$ws.Range("A89").Value = "Croacia"
$ws.Range("B89").Value = 13368
$ws.Range("C89").Value = 261
$ws.Range("D89").Value = 10721
$ws.Range("E89").Value = 2429
$ws.Range("G89").Value = 7
$ws.Range("H89").Value = 218

$ws.Range("A90").Value = "Zambia"
$ws.Range("B90").Value = 13323
$ws.Range("C90").Value = 0
$ws.Range("D90").Value = 11899
$ws.Range("E90").Value = 1118
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 306

# Estonia (row 141)
$ws.Range("B141").Value = 2655
$ws.Range("C141").Value = 23
$ws.Range("D141").Value = 2252
$ws.Range("E141").Value = 339
